# Delete the row containing "RIVO 320MG 20*10 TABS" (row 19), which
# removes that item from the shortage list and shifts all subsequent
# rows (items, totals and footer) up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Delete()
